$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Cuaderno" sheet as a copy of "Tareas", placed right
#    after it, then adapt it into the diagnostic-exam / notebook tracker.
# ---------------------------------------------------------------------------
$tareas = $wb.Worksheets.Item("Tareas")
$tareas.Copy($null, $tareas)
$cuaderno = $wb.Worksheets.Item("Tareas (2)")
$cuaderno.Name = "Cuaderno"

# Insert a brand-new row above the old header row (old row 3 -> row 4, etc.)
$cuaderno.Rows("3").Insert()

# Old row 3 (the "Tarea"/1..15/"Promedio" row) is now row 4 and keeps the
# blue fill; bring in its formatting for the new row 3 first so D3/E3:S3
# start from the same base (bold, blue fill, bordered, centered) seen
# throughout the workbook's other header rows.
$cuaderno.Range("D4").Copy()
$cuaderno.Range("D3").PasteSpecial(-4122)
$cuaderno.Range("E4").Copy()
$cuaderno.Range("E3:S3").PasteSpecial(-4122)
$cuaderno.Range("E4").Copy()
$cuaderno.Range("B3:C3").PasteSpecial(-4122)
$cuaderno.Range("E4").Copy()
$cuaderno.Range("T3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 3 content: "Cuaderno" / "Examen diagnóstico" (tall, wrapped,
# rotated header like a scored-items grid).
$cuaderno.Range("E3").Value = "Examen diagnóstico"
$cuaderno.Range("D3").Value = "Cuaderno"

$cuaderno.Range("D3:T3").VerticalAlignment = -4108
$cuaderno.Range("D3").HorizontalAlignment = -4108
$cuaderno.Range("E3:S3").HorizontalAlignment = -4108
$cuaderno.Range("E3:S3").Orientation = 90
$cuaderno.Range("E3:S3").WrapText = $true
$cuaderno.Rows(3).RowHeight = 123.75

# Row 4 (old row 3) loses its blue fill and becomes plain/white, but keeps
# borders and gains vertical centering.
$cuaderno.Range("B4:T4").Interior.Pattern = -4142
$cuaderno.Range("B4:T4").VerticalAlignment = -4108
$cuaderno.Range("B4:T4").HorizontalAlignment = -4108
$cuaderno.Range("D4").Value = "Actividad"

# ---------------------------------------------------------------------------
# 2. Make "Cuaderno" the active/selected sheet (mirrors the workbook-level
#    activeTab change in the diff) and drop the old tab selection flag from
#    "Calificaciones".
# ---------------------------------------------------------------------------
$cuaderno.Select() | Out-Null
$cuaderno.Range("B2").Select() | Out-Null

Write-Host "done"
